$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update I2 value
$ws.Range("I2").Value = 11.73

# Add new I column values for rows 3-7
$ws.Range("I3").Value = 11.85
$ws.Range("I4").Value = 12.71
$ws.Range("I5").Value = 11.88
$ws.Range("I6").Value = 12.65
$ws.Range("I7").Value = 12.26

# Update the selection (active cell) to D9
$ws.Range("D9").Select()

# Update tab color to opaque white (RGB 255,255,255 -> 255 + 255*256 + 255*65536)
$ws.Tab.Color = 16777215
